$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, using the same style as the other headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H11 with 0
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

# H12 gets value 1
$ws.Cells.Item(12, 8).Value = 1
